$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Remove the rows that were dropped (GRANITO, LED, PANELES PU, PANELES PVC) ---
# Delete bottom-up so earlier row numbers stay valid while deleting.
$ws.Range("A14").EntireRow.Delete()
$ws.Range("A13").EntireRow.Delete()
$ws.Range("A9").EntireRow.Delete()
$ws.Range("A5").EntireRow.Delete()

# --- Refresh the remaining data rows (2-14) with their updated figures ---
$ws.Cells.Item(2,3).Value = 6197.58402943659
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 6197.58402943659
$ws.Cells.Item(2,6).Value = 0

$ws.Cells.Item(3,3).Value = 17669.1470988183
$ws.Cells.Item(3,4).Value = 1102.46
$ws.Cells.Item(3,5).Value = 16566.6870988183
$ws.Cells.Item(3,6).Value = 0.06239463590598167

$ws.Cells.Item(4,3).Value = 1043.22288526528
$ws.Cells.Item(4,4).Value = 189.6
$ws.Cells.Item(4,5).Value = 853.6228852652799
$ws.Cells.Item(4,6).Value = 0.1817444792267827

$ws.Cells.Item(5,3).Value = 150
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 150
$ws.Cells.Item(5,6).Value = 0

$ws.Cells.Item(6,3).Value = 2907.58368146026
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 2907.58368146026
$ws.Cells.Item(6,6).Value = 0

$ws.Cells.Item(7,3).Value = 886.711016287574
$ws.Cells.Item(7,4).Value = 313.2
$ws.Cells.Item(7,5).Value = 573.511016287574
$ws.Cells.Item(7,6).Value = 0.3532154154476236

$ws.Cells.Item(8,3).Value = 1346.40488751609
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 1346.40488751609
$ws.Cells.Item(8,6).Value = 0

$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 0

$ws.Cells.Item(10,3).Value = 3881.07983534392
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 3881.07983534392
$ws.Cells.Item(10,6).Value = 0

$ws.Cells.Item(11,3).Value = 17831.4143984654
$ws.Cells.Item(11,4).Value = 2486.85
$ws.Cells.Item(11,5).Value = 15344.5643984654
$ws.Cells.Item(11,6).Value = 0.1394645396281084

$ws.Cells.Item(12,3).Value = 61863.7203947566
$ws.Cells.Item(12,4).Value = 4175.49
$ws.Cells.Item(12,5).Value = 57688.2303947566
$ws.Cells.Item(12,6).Value = 0.06749497077375745

$ws.Cells.Item(13,3).Value = 440.653177778119
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 440.653177778119
$ws.Cells.Item(13,6).Value = 0

$ws.Cells.Item(14,3).Value = 7837.31410570622
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 7837.31410570622
$ws.Cells.Item(14,6).Value = 0

# --- TOTAL row (15) ---
$ws.Cells.Item(15,3).Value = 122054.8355108344
$ws.Cells.Item(15,4).Value = 8267.599999999999
$ws.Cells.Item(15,5).Value = 113787.2355108344
$ws.Cells.Item(15,6).Value = 0.06773676737507146

# --- Column widths: D -> 13, E -> 23, F -> 25 (stored width = ColumnWidth + 5/6) ---
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 22.166666666666668
$ws.Columns.Item(6).ColumnWidth = 24.166666666666668
